# v1.1 verify the status after modification
# verify the status and closed after modification done on the navigation test cases

$wb = $excel.ActiveWorkbook

# --- Sheet: LH_TC_NAVIGATION _REVIEWS ---
# Update "Reviewer verification" (column J) from "Open" to "Closed" for all review rows (2-6)
$wsReviews = $wb.Worksheets.Item("LH_TC_NAVIGATION _REVIEWS")
$wsReviews.Range("J2:J6").Value = "Closed"

# --- Sheet: Version History ---
# Add a new version history row documenting this update
$wsHistory = $wb.Worksheets.Item("Version History")
$wsHistory.Range("A3").Value = "v1.1"
$wsHistory.Range("B3").Value = "Ahmed Abuzaid"
$wsHistory.Range("C3").Value = "verify the status after modification done on the test cases"
$wsHistory.Range("D3").Formula = "=TODAY()"
$wsHistory.Rows.Item(3).RowHeight = 30

# --- Selections / active sheet, to mirror the saved UI state ---
$wsReviews.Range("J6").Select()
$wsHistory.Select()
$wsHistory.Range("C3").Select()
